# Append 5 new daily rows (234-238) to Sheet1, continuing the existing
# "somma mobile 7gg." series through 2021-04-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data: date (Excel serial), B (nuovi pos.), C (somma mobile 7gg.), D (per 100mila ab.)
$newRows = @(
    @(44308, 0, 7, 106.6098081023454),
    @(44309, 2, 7, 106.6098081023454),
    @(44310, 2, 8, 121.8397806883948),
    @(44311, 0, 5, 76.14986293024673),
    @(44312, 1, 6, 91.37983551629607)
)

$startRow = 234
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    # Copy formatting (style) from the cell directly above (column A uses
    # the date-format / bordered style "s=2") before writing the new value.
    $ws.Range("A" + ($r - 1)).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}

$excel.CutCopyMode = 0
